$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(533).Insert()

$ws.Range("A533").Value = 4
$ws.Range("B533").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C533").Value = "Los Lagos"
$ws.Range("D533").Value = 45258
$ws.Range("E533").Value = 10
$ws.Range("F533").Value = 100112008
$ws.Range("G533").Value = "Coliflor"
$ws.Range("H533").Value = "Sin especificar"
$ws.Range("I533").Value = "Primera"
$ws.Range("J533").Value = 1500
$ws.Range("K533").Value = 1600
$ws.Range("L533").Value = 1600
$ws.Range("M533").Value = 1600
$ws.Range("N533").Value = "`$/unidad"
$ws.Range("O533").Value = "Región Metropolitana"
$ws.Range("P533").Value = 1600
$ws.Range("Q533").Value = 1
$ws.Range("R533").Value = "Hortaliza"
